# Applies the "Updated cryptos list" data refresh to Sheet1.
# For every touched cell we pin NumberFormat to "@" (Text) before writing the
# new value so price strings such as "229.20" / "0.0610" / "1.00" keep their
# original trailing zeros instead of being auto-coerced to numbers by COM.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.596.85"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.41%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.632.64"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.60%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.42"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.69%  "

# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.85%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.07"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.09%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.263"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0610"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.32%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0861"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.40%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.864.26"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.57%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.638.93"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.08%  "

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.50%  "

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.42%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.07"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.593.66"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.36%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.20"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.75%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.54"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.70%  "

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.17%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.63"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +6.44%  "

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.17%  "

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.49%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.13"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.79%  "

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.12%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.77%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.62"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.48%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.09%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.70%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.92%  "

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.36%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.462.13"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.70%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.54"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.44%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.31"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.88%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.878"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.37%  "

# Row 38
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.925"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.46%  "

# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.558"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.66%  "

# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0167"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.71%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "68.93"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.00%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.08%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.72%  "

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.03%  "

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.29%  "

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.25%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.773.95"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.71%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.74"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.97%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.37"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.17%  "

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.72%  "

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.14%  "
